$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (e.g. trailing zeros,
# multi-dot "thousands" separators, padded percentages) by forcing Text
# number format before assigning the literal string value.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.826.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.113.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.90'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.109.12'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.47'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.96'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.629.98'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.783.33'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.112.10'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.26'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '476.17'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.713'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.95'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.28%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.32'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.51%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.09%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.85'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.56'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -7.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '47.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.21'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.05'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.829.45'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '384.51'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.39%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '135.43'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.75'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.44%  '
